# Update "想去人数" (want-to-go count, column F) values for the first five
# event rows (rows 2-6) on both the "展览" sheet and the "全部类型" sheet,
# matching the regenerated site data.

$wb = $excel.ActiveWorkbook

$newValues = @{
    2 = 3438
    3 = 25
    4 = 68
    5 = 1714
    6 = 90
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $newValues[$row]
    }
}
